$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.407.72"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "'3.840.98"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'602.07"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'169.26"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "'3.843.39"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "'0.0000268"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").Value = "'37.14"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'4.489.18"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "'3.851.20"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "'68.470.54"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "'18.53"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'470.36"
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").Value = "'0.735"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'0.0000161"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").Value = "'83.32"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'3.994.28"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").Value = "'31.58"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("D35").Value = "'9.38"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").Value = "'3.809.65"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("E38").Value = "  +10.35%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "'5.95"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("E44").Value = "  -4.91%  "
$ws.Range("D45").Value = "'8.73"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D47").Value = "'414.53"
$ws.Range("E47").Value = "  -5.23%  "
$ws.Range("D48").Value = "'47.15"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  +4.70%  "
$ws.Range("D50").Value = "'0.0361"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "'141.57"
$ws.Range("E51").Value = "  -1.70%  "
